# Generate Report for Handoff
#
# A new handoff round has completed for the "f2c8c8bd-1352-4160-bf64-93a10ac92f7b.md"
# file (the last row on every sheet). Refresh the handoff timestamps so the
# report reflects the freshly generated xliff files:
#   - Overview sheet: "Latest HO Xliff Generate Date" for that file
#   - zh-cn sheet:     "Latest Handoff Datetime" for that file
#   - de-de sheet:     "Latest Handoff Datetime" for that file

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-14 16:56:21"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-14 16:56:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-14 16:56:21"
